$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  B = "RPA Uipath Developer"; C = "https://www.naukri.com/job-listings-rpa-uipath-developer-semperfi-solutions-bengaluru-13-to-15-years-190524902070" },
    @{ Row = 3;  B = "RPA Developer"; C = "https://www.naukri.com/job-listings-rpa-developer-percipere-mumbai-0-to-1-years-190524902221" },
    @{ Row = 4;  B = "Application Automation Engineer"; C = "https://www.naukri.com/job-listings-application-automation-engineer-accenture-solutions-pvt-ltd-pune-3-to-5-years-190524901638" },
    @{ Row = 5;  B = "Service Management New Associate - Automation and Analytics"; C = "https://www.naukri.com/job-listings-service-management-new-associate-automation-and-analytics-accenture-solutions-pvt-ltd-bengaluru-0-to-1-years-190524901762" },
    @{ Row = 6;  B = "Trust & Safety New Associate"; C = "https://www.naukri.com/job-listings-trust-safety-new-associate-accenture-solutions-pvt-ltd-gurugram-0-to-1-years-180524908538" },
    @{ Row = 7;  B = "BPM Consultant"; C = "https://www.naukri.com/job-listings-bpm-consultant-percipere-mumbai-7-to-11-years-180524906883" },
    @{ Row = 8;  B = "Insurance Operations Manager"; C = "https://www.naukri.com/job-listings-insurance-operations-manager-accenture-solutions-pvt-ltd-hyderabad-16-to-25-years-180524903816" },
    @{ Row = 9;  B = "Application Developer"; C = "https://www.naukri.com/job-listings-application-developer-accenture-solutions-pvt-ltd-mumbai-3-to-5-years-190524901581" },
    @{ Row = 10; B = "Service Management Senior Analyst"; C = "https://www.naukri.com/job-listings-service-management-senior-analyst-accenture-solutions-pvt-ltd-chennai-5-to-8-years-180524908644" },
    @{ Row = 11; B = "Service Management Analyst"; C = "https://www.naukri.com/job-listings-service-management-analyst-accenture-solutions-pvt-ltd-chennai-3-to-5-years-180524904105" },
    @{ Row = 12; B = "Service Management Associate"; C = "https://www.naukri.com/job-listings-service-management-associate-accenture-solutions-pvt-ltd-bengaluru-1-to-3-years-190524903201" },
    @{ Row = 13; B = "Application Developer"; C = "https://www.naukri.com/job-listings-application-developer-accenture-solutions-pvt-ltd-bengaluru-3-to-5-years-180524906353" },
    @{ Row = 14; B = "Technology Architect"; C = "https://www.naukri.com/job-listings-technology-architect-accenture-solutions-pvt-ltd-bengaluru-12-to-16-years-190524903287" },
    @{ Row = 15; B = "Business Analyst"; C = "https://www.naukri.com/job-listings-business-analyst-accenture-solutions-pvt-ltd-pune-5-to-9-years-180524904126" },
    @{ Row = 16; B = "Underwriting Specialist"; C = "https://www.naukri.com/job-listings-underwriting-specialist-accenture-solutions-pvt-ltd-gurugram-7-to-11-years-180524903796" },
    @{ Row = 17; B = "SW/App/Cloud Tech Support Analyst"; C = "https://www.naukri.com/job-listings-sw-app-cloud-tech-support-analyst-accenture-solutions-pvt-ltd-hyderabad-3-to-5-years-180524907692" },
    @{ Row = 18; B = "Service Management Associate"; C = "https://www.naukri.com/job-listings-service-management-associate-accenture-solutions-pvt-ltd-chennai-1-to-3-years-180524903797" },
    @{ Row = 19; B = "Data Analyst - IO - Client Reporting - Investment Data Management"; C = "https://www.naukri.com/job-listings-data-analyst-io-client-reporting-investment-data-management-m-amp-amp-amp-g-plc-mumbai-1-to-4-years-180524500057" },
    @{ Row = 20; B = "Analyst - Revenue and Fixed assets profile"; C = "https://www.naukri.com/job-listings-analyst-revenue-and-fixed-assets-profile-m-amp-amp-amp-g-plc-mumbai-1-to-3-years-180524500056" },
    @{ Row = 21; B = "Assistant Manager - Actuarial"; C = "https://www.naukri.com/job-listings-assistant-manager-actuarial-m-amp-amp-amp-g-plc-mumbai-3-to-5-years-180524500055" }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}
